$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.133.05'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.653.67'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.51'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5293'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2610'
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06335'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.41'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07759'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.675.81'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.487'
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5464'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8132'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.29'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.126.51'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.548'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.04'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.995'
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '140.26'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.265'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.432'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05941'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.511'
$ws.Range('E31').Value = '  -5.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.240'
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.551'
$ws.Range('E33').Value = '  -4.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.413'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9456'
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5639'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01608'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.859'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8455'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.009.00'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.84'
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.799.69'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '56.79'
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈105'
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.005'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.471'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05152'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.744'
$ws.Range('E51').Value = '  -4.26%  '
